$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.794.28'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').Value = '2.093.93'
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.11'
$ws.Range('E5').Value = '  -1.09%  '
$ws.Range('E6').Value = '  -2.36%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.43'
$ws.Range('E8').Value = '  -5.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '59.66'
$ws.Range('E9').Value = '  -1.13%  '
$ws.Range('E10').Value = '  -4.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0768'
$ws.Range('E11').Value = '  -2.45%  '
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.00'
$ws.Range('E13').Value = '  -5.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.882'
$ws.Range('E14').Value = '  +5.46%  '
$ws.Range('D15').Value = '2.396.95'
$ws.Range('E15').Value = '  +1.84%  '
$ws.Range('E16').Value = '  -3.77%  '
$ws.Range('D17').Value = '2.090.35'
$ws.Range('E17').Value = '  +1.96%  '
$ws.Range('D18').Value = '36.751.28'
$ws.Range('E18').Value = '  -1.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.44'
$ws.Range('E19').Value = '  -3.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.24'
$ws.Range('E20').Value = '  -2.60%  '
$ws.Range('D21').Value = '0.0₃0878'
$ws.Range('E21').Value = '  -2.59%  '
$ws.Range('E22').Value = '  +1.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.74'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('E25').Value = '  -2.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.90'
$ws.Range('E26').Value = '  +5.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.17'
$ws.Range('E27').Value = '  -1.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '168.20'
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.89'
$ws.Range('E29').Value = '  +3.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.39'
$ws.Range('E30').Value = '  +11.29%  '
$ws.Range('E31').Value = '  -0.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.19'
$ws.Range('E32').Value = '  +5.99%  '
$ws.Range('E33').Value = '  +4.03%  '
$ws.Range('E34').Value = '  -1.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.41'
$ws.Range('E35').Value = '  +5.37%  '
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('E37').Value = '  +4.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0838'
$ws.Range('E38').Value = '  -6.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.28'
$ws.Range('E39').Value = '  -4.05%  '
$ws.Range('E40').Value = '  +1.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.91'
$ws.Range('E41').Value = '  -5.52%  '
$ws.Range('E42').Value = '  -0.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0958'
$ws.Range('E43').Value = '  -8.16%  '
$ws.Range('E44').Value = '  -8.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '96.71'
$ws.Range('E45').Value = '  +0.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.37'
$ws.Range('E46').Value = '  -5.54%  '
$ws.Range('D47').Value = '1.347.25'
$ws.Range('E47').Value = '  +4.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.42'
$ws.Range('E48').Value = '  -1.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.10'
$ws.Range('E49').Value = '  +3.60%  '
$ws.Range('E50').Value = '  -1.05%  '
$ws.Range('D51').Value = '2.279.81'
$ws.Range('E51').Value = '  +1.65%  '
